$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$casesQuery = @'
MATCH (c:case)
    WHERE c.gender='MALE'
OPTIONAL MATCH (c)-[:of_arm]->(a:arm)-[:of_trial]->(ct:clinical_trial)
OPTIONAL MATCH (f:file)-[*]->(c)
RETURN DISTINCT
    c.case_id AS `Case ID`,
     ct.clinical_trial_designation AS `Trial Code`,
     a.arm_id AS Arm,
      a.arm_drug AS `Arm Treatment`,
c.disease AS Diagnosis,
  c.gender AS Gender,
    c.race AS Race,
    c.ethnicity AS Ethnicity
 
'@

$statQuery = @'
MATCH (s:specimen)-->(c:case)-->(:arm)-->(ct:clinical_trial)
    WHERE c.gender = "MALE"
OPTIONAL MATCH (f:file)-->(:sequencing_assay)-->(:nucleic_acid)-->(s)
RETURN 
     COUNT(DISTINCT ct.clinical_trial_designation) AS Trials,
     COUNT(DISTINCT c.case_id) AS Cases,
      COUNT(DISTINCT f) AS Files
'@

$filesQuery = @'

MATCH (c:case)
 WHERE c.gender='MALE'
   MATCH (f:file)
      OPTIONAL MATCH (f)-[*]->(a:arm)-[:of_trial]->(ct:clinical_trial)
       MATCH (f)-[*]->(c)
      OPTIONAL MATCH (f)-->(parent)
      WITH
        f, parent, c, a, ct, 
        ['Bytes', 'KB', 'MB', 'GB', 'TB'] AS units,
        toInteger(floor(log(f.file_size)/log(1024))) as i,
        2 as precision
WITH
        f, parent, c, a, ct,
        f.file_size /(1024^i) AS value, 
        10^precision AS factor,
        units[i] as unit
WITH    
        f, parent, c, a, ct, unit,
        round(factor * value)/factor AS size
      RETURN DISTINCT 
       f.file_name AS `File Name`,
       head(labels(parent)) as Association,
       f.file_description AS Description,
       f.file_format AS `File Format`,
         CASE size % 1 WHEN 0 THEN apoc.convert.toInteger(size)+' ' +unit ELSE size+' ' +unit END AS Size,
       ct.clinical_trial_designation AS `Trial Code`,
       a.arm_id AS Arm,
       c.case_id AS `Case ID`
        
'@

$ws.Range("B2").Value = $casesQuery
$ws.Range("C2").Value = $statQuery
$ws.Range("B3").Value = $filesQuery
$ws.Range("C3").Value = $statQuery

$ws.Rows.Item(2).RowHeight = 210
$ws.Rows.Item(3).RowHeight = 409.5

$ws.Range("D3").Select()
$excel.ActiveWindow.ScrollRow = 3
$excel.ActiveWindow.ScrollColumn = 1

Write-Host "done"
